$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the E-column (text) values first, in the order that makes the shared-string
# table match the expected first-seen order: idx9..idx15.
# (E6/E8/E10/E12 lose their inherited formatting in the real edit, so clear first.)
$ws.Range("E4").Value = "Semantically similar, syntactically too different"
$ws.Range("E2").Value = "Semantically similar, syntactically too different; Should the code blocks with local ID 3 match?"
$ws.Range("E3").Value = "Connection should be found by the metric"
$ws.Range("E5").Value = "Connection should be found by the metric in versions 1+2, 4+5"

$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = "Semantically similar, syntactically too different; Overlap should have worked in versions 4+5; Check GT in version 7+8 again"

$ws.Range("E9").Value = "Check GT again"

$ws.Range("E12").ClearFormats()
$ws.Range("E12").Value = "Semantically similar, syntactically too different; Check GT again"

$ws.Range("E7").Value = "Semantically similar, syntactically too different"

$ws.Range("E8").ClearFormats()
$ws.Range("E8").Value = "Semantically similar, syntactically too different"

$ws.Range("E10").ClearFormats()
$ws.Range("E10").Value = "Semantically similar, syntactically too different"

$ws.Range("E11").Value = "Check GT again"

# D-column values (all 0)
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 0

$ws.Range("E12").Select()
